$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Role" column (C), shifting
# Role/Gender/Age one column to the right (now D/E/F).
$ws.Columns("C:C").Insert()

# New "Password" column: header + a "pwd" placeholder value for every
# staff member (CSV-import style placeholder for the login password).
$ws.Range("C1").Value = "Password"
$ws.Range("C2").Value = "pwd"
$ws.Range("C3").Value = "pwd"
$ws.Range("C4").Value = "pwd"
$ws.Range("C5").Value = "pwd"

# Match the bold/centered/bordered header formatting used by the other
# header cells (copy formats only, so the new values are untouched).
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-balance the column widths now that there's an extra column between
# "Name" and "Role".
$ws.Columns("B:C").ColumnWidth = 23.83
$ws.Columns("D:D").ColumnWidth = 29.66

# Move the active selection.
$ws.Range("C9").Select()
